# Append the "Post 59" row (While Loop | Shell Scripting) to the blog log
# table on Sheet1, growing the table/autofilter and sheet dimension to
# cover the new row, and updating the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 68 is the last existing data row of the table; clone its
# formatting onto the new row 69 before filling in the new values so the
# new row looks consistent with the rest of the table (plain style for
# S.No/Title/Date, Hyperlink style for the two link columns).
$ws.Range("B68:F68").Copy()
$ws.Range("B69:F69").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B69").Value2 = 59
$ws.Range("C69").Value2 = "While Loop | Shell Scripting"
$ws.Range("D69").Value2 = 44178
$ws.Range("E69").Value2 = "https://programmingport.hashnode.dev/while-loop-or-shell-scripting"
$ws.Range("F69").Value2 = "https://dev.to/rahulmishra05/while-loop-shell-scripting-5f9a"

# Grow the table (Table2) and its autofilter to include the new row.
$lo = $ws.ListObjects.Item("Table2")
$lo.Resize($ws.Range("B10:F69"))

# Update the selection to match the newly added last cell.
$null = $ws.Range("F69").Select()
